# Re-sort the per-source statistics rows (B2:X8) so that each database
# name lines up with its correct precision/recall/fmeasure figures,
# matching the refreshed dataframe dump (see commit "dataframe, chart and
# styles to visualize"). Row 1 (headers) and row 9 (the union totals row)
# are unaffected by the re-sort.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "wiley"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 295
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 295
$ws.Range("P2").Value = 294
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0

# Row 4
$ws.Range("B4").Value = "springer"
$ws.Range("C4").Value = 0.8099999999999999
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 1.39
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 124
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 124
$ws.Range("P4").Value = 115
$ws.Range("R4").Value = 9
$ws.Range("S4").Value = 0.8099999999999999
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = 1.39

# Row 5
$ws.Range("B5").Value = "scopus"
$ws.Range("C5").Value = 9.09
$ws.Range("D5").Value = 35
$ws.Range("E5").Value = 14.43
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 77
$ws.Range("M5").Value = 7
$ws.Range("N5").Value = 77
$ws.Range("P5").Value = 68
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 9
$ws.Range("S5").Value = 9.09
$ws.Range("T5").Value = 35
$ws.Range("U5").Value = 14.43
$ws.Range("V5").Value = 0.02597402597402598
$ws.Range("X5").Value = 0.04705882352941176

# Row 6
$ws.Range("B6").Value = "sciencedirect"
$ws.Range("C6").Value = 2.41
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 4.46
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 249
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 249
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 245
$ws.Range("R6").Value = 4
$ws.Range("S6").Value = 2.41
$ws.Range("T6").Value = 30
$ws.Range("U6").Value = 4.46
$ws.Range("V6").Value = 0.008032128514056224
$ws.Range("W6").Value = 0.25
$ws.Range("X6").Value = 0.01556420233463035

# Row 7
$ws.Range("B7").Value = "elcompendex"
$ws.Range("C7").Value = 25
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 14.29
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 8
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 8
$ws.Range("P7").Value = 1
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = 10
$ws.Range("U7").Value = 14.29

# Row 8
$ws.Range("B8").Value = "acm"
$ws.Range("C8").Value = 2.38
$ws.Range("D8").Value = 25
$ws.Range("E8").Value = 4.35
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 210
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 210
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 210
$ws.Range("Q8").Value = 5
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 2.38
$ws.Range("T8").Value = 25
$ws.Range("U8").Value = 4.35
$ws.Range("V8").Value = 0.01428571428571429
$ws.Range("W8").Value = 0.375
$ws.Range("X8").Value = 0.02752293577981652
